$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff shows cell C10 (value "R20" row, Integer min column) changing
# from 18 to 1.
$ws.Range("C10").Value = 1
